$d = $word.ActiveDocument

# 1. "See more on stakeholders document and diagram." ->
#    "See more on Stakeholder's document and diagram."
$d.Content.Find.Execute(
    "See more on stakeholders document and diagram.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "See more on Stakeholder's document and diagram.", 2)

# 2. Relocate the "_GoBack" bookmark: it used to sit right after the
#    paragraph above ("...diagram."); it now belongs inside the word
#    "interact" in the later "Built-in messaging system" bullet, right
#    between "intera" and "ct".
$rng = $d.Content
$rng.Find.Execute("to intera", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$target = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $target)
